$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.076.09'
$ws.Range('E2').Value = '  -0.57%  '
$ws.Range('D3').Value = '1.652.15'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = "'217.15"
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').Value = "'0.5279"
$ws.Range('E6').Value = '  +0.98%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = "'0.2599"
$ws.Range('E8').Value = '  -1.63%  '
$ws.Range('D9').Value = "'0.06314"
$ws.Range('E9').Value = '  +0.62%  '
$ws.Range('D10').Value = "'20.33"
$ws.Range('E10').Value = '  -2.24%  '
$ws.Range('D11').Value = "'0.07788"
$ws.Range('E11').Value = '  +0.35%  '
$ws.Range('D12').Value = "'4.516"
$ws.Range('E12').Value = '  +0.98%  '
$ws.Range('D13').Value = '1.648.39'
$ws.Range('E13').Value = '  -3.39%  '
$ws.Range('D14').Value = '1.879.26'
$ws.Range('E14').Value = '  -0.64%  '
$ws.Range('D15').Value = "'0.5473"
$ws.Range('E15').Value = '  +0.39%  '
$ws.Range('D16').Value = '0.0₅8192'
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('D17').Value = "'65.32"
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('D18').Value = '26.076.42'
$ws.Range('E18').Value = '  -0.59%  '
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D20').Value = "'4.579"
$ws.Range('E20').Value = '  -0.57%  '
$ws.Range('D21').Value = "'190.59"
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('D22').Value = "'10.06"
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('D23').Value = "'6.012"
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').Value = "'143.93"
$ws.Range('E25').Value = '  +3.66%  '
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('D27').Value = "'7.207"
$ws.Range('E27').Value = '  -0.67%  '
$ws.Range('D28').Value = "'15.98"
$ws.Range('E28').Value = '  -1.05%  '
$ws.Range('D29').Value = "'1.457"
$ws.Range('E29').Value = '  +3.07%  '
$ws.Range('D30').Value = "'0.05796"
$ws.Range('E30').Value = '  -2.83%  '
$ws.Range('D31').Value = "'1.271"
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('D32').Value = "'3.541"
$ws.Range('E32').Value = '  +0.10%  '
$ws.Range('D33').Value = "'3.261"
$ws.Range('E33').Value = '  -0.14%  '
$ws.Range('D34').Value = "'1.591"
$ws.Range('E34').Value = '  +0.83%  '
$ws.Range('D36').Value = "'2.412"
$ws.Range('E36').Value = '  -0.34%  '
$ws.Range('D37').Value = "'0.9425"
$ws.Range('E37').Value = '  -1.71%  '
$ws.Range('E38').Value = '  +1.17%  '
$ws.Range('D40').Value = "'0.8495"
$ws.Range('E40').Value = '  -0.27%  '
$ws.Range('D41').Value = "'104.37"
$ws.Range('E41').Value = '  +3.78%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('E43').Value = '  -4.30%  '
$ws.Range('D44').Value = '1.030.60'
$ws.Range('E44').Value = '  +2.51%  '
$ws.Range('D45').Value = '1.793.14'
$ws.Range('E45').Value = '  -0.70%  '
$ws.Range('D46').Value = "'56.93"
$ws.Range('E46').Value = '  +0.46%  '
$ws.Range('E47').Value = '  -0.16%  '
$ws.Range('D48').Value = "'0.4328"
$ws.Range('E48').Value = '  -0.42%  '
$ws.Range('D49').Value = "'7.846"
$ws.Range('E49').Value = '  -1.92%  '
$ws.Range('E50').Value = '  -0.18%  '
$ws.Range('D51').Value = "'1.446"
$ws.Range('E51').Value = '  -1.37%  '
